$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.663.49'
$ws.Range("E2").Value = '  +1.22%  '

$ws.Range("D3").Value = '1.868.45'
$ws.Range("E3").Value = '  +0.38%  '

$ws.Range("E4").Value = '  +0.29%  '

$ws.Range("D5").Value = '''331.92'
$ws.Range("E5").Value = '  +2.58%  '

$ws.Range("E6").Value = '  +0.23%  '

$ws.Range("D7").Value = '''0.4693'
$ws.Range("E7").Value = '  +3.87%  '

$ws.Range("D8").Value = '''0.3938'
$ws.Range("E8").Value = '  +1.85%  '

$ws.Range("D9").Value = '''47.79'
$ws.Range("E9").Value = '  -0.60%  '

$ws.Range("D10").Value = '''0.08042'
$ws.Range("E10").Value = '  +1.74%  '

$ws.Range("E11").Value = '  +0.15%  '

$ws.Range("D12").Value = '''21.86'
$ws.Range("E12").Value = '  +1.96%  '

$ws.Range("D13").Value = '1.873.70'
$ws.Range("E13").Value = '  +0.83%  '

$ws.Range("D14").Value = '''5.953'
$ws.Range("E14").Value = '  +1.00%  '

$ws.Range("D15").Value = '''7.136'
$ws.Range("E15").Value = '  -0.23%  '

$ws.Range("E16").Value = '  +0.26%  '

$ws.Range("D17").Value = '''0.00001045'
$ws.Range("E17").Value = '  +0.96%  '

$ws.Range("D18").Value = '''86.61'
$ws.Range("E18").Value = '  +0.87%  '

$ws.Range("D19").Value = '''0.06637'
$ws.Range("E19").Value = '  +1.71%  '

$ws.Range("D20").Value = '''17.09'
$ws.Range("E20").Value = '  -0.63%  '

$ws.Range("D21").Value = '''1.002'
$ws.Range("E21").Value = '  +0.26%  '

$ws.Range("D22").Value = '27.669.38'
$ws.Range("E22").Value = '  +1.23%  '

$ws.Range("D23").Value = '''5.507'
$ws.Range("E23").Value = '  -0.21%  '

$ws.Range("D24").Value = '''10.97'
$ws.Range("E24").Value = '  +1.29%  '

$ws.Range("E25").Value = '  +1.86%  '

$ws.Range("D26").Value = '2.098.79'
$ws.Range("E26").Value = '  +0.75%  '

$ws.Range("D27").Value = '''158.65'
$ws.Range("E27").Value = '  +4.07%  '

$ws.Range("D28").Value = '''20.18'
$ws.Range("E28").Value = '  +2.21%  '

$ws.Range("D29").Value = '''2.089'
$ws.Range("E29").Value = '  +1.15%  '

$ws.Range("D30").Value = '''5.551'
$ws.Range("E30").Value = '  +0.86%  '

$ws.Range("D31").Value = '''122.25'
$ws.Range("E31").Value = '  +1.43%  '

$ws.Range("D32").Value = '''0.9696'
$ws.Range("E32").Value = '  +3.46%  '

$ws.Range("D33").Value = '''0.09507'
$ws.Range("E33").Value = '  +2.09%  '

$ws.Range("D34").Value = '''1.448'
$ws.Range("E34").Value = '  -3.27%  '

$ws.Range("E35").Value = '  -0.39%  '

$ws.Range("D36").Value = '''5.318'
$ws.Range("E36").Value = '  +0.66%  '

$ws.Range("D37").Value = '''0.02259'
$ws.Range("E37").Value = '  +1.01%  '

$ws.Range("D38").Value = '''0.06095'
$ws.Range("E38").Value = '  +1.47%  '

$ws.Range("D39").Value = '''1.224'
$ws.Range("E39").Value = '  -0.03%  '

$ws.Range("D40").Value = '''8.129'
$ws.Range("E40").Value = '  -1.49%  '

$ws.Range("B41").Value = 'Frax'
$ws.Range("C41").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D41").Value = '''1.002'
$ws.Range("E41").Value = '  +0.21%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '''0.5993'
$ws.Range("E42").Value = '  +1.29%  '

$ws.Range("D43").Value = '''0.1887'
$ws.Range("E43").Value = '  -0.07%  '

$ws.Range("D44").Value = '''10.21'
$ws.Range("E44").Value = '  +0.30%  '

$ws.Range("E45").Value = '  -1.04%  '

$ws.Range("E46").Value = '  +0.70%  '

$ws.Range("D47").Value = '''12.21'
$ws.Range("E47").Value = '  +2.16%  '

$ws.Range("D48").Value = '''3.389'
$ws.Range("E48").Value = '  +0.60%  '

$ws.Range("D49").Value = '''1.932'
$ws.Range("E49").Value = '  +0.26%  '

$ws.Range("D50").Value = '''0.06850'
$ws.Range("E50").Value = '  +0.80%  '

$ws.Range("D51").Value = '''114.58'
$ws.Range("E51").Value = '  +6.19%  '
